$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 6 - Воробьев Иван Александрович: B6 0 -> 5
$ws.Range("B6").Value = 5

# Row 7 - Глазков Илья Максимович: B7 (empty) -> 5
$ws.Range("B7").Value = 5

# Row 8 - Губеева София Юрьевна: add note in M8
$ws.Range("M8").Value = "3й защищен"

# Row 11 - Капелина Анна Ивановна: B11 (empty) -> 5
$ws.Range("B11").Value = 5

# Row 15 - Лунин Владимир Анатольевич: D15 0 -> 5 (J15/K15 formulas recalc automatically)
$ws.Range("D15").Value = 5

# Row 16 - Оганезов Михаил Алексеевич: B16 0 -> 5
$ws.Range("B16").Value = 5

# Row 17 - Попова Полина Владиславовна: B17 (empty) -> 5
$ws.Range("B17").Value = 5

# Row 24 - Титова Надежда Алексеевна: B24 0 -> 5
$ws.Range("B24").Value = 5

# Row 25 - Хрищанович Полина Чеславовна: B25 (empty) -> 5
$ws.Range("B25").Value = 5
